$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in WEEK 14 (row 13) ---
$ws.Range("A13").Value = "WEEK 14"
$ws.Range("B13").Value = "N/A"
$ws.Range("C13").Value = "N/A"
$ws.Range("D13").Value = "N/A"
$ws.Range("E13").Value = "N/A"
$ws.Range("B13:E13").WrapText = $true

# --- Fill in WEEK 15 (row 14) ---
$ws.Range("A14").Value = "WEEK 15"
$ws.Range("B14").Value = "Shane, Jaylee, Cameron, Ryan, Angelo, Chandler"
$ws.Range("C14").Value = "N/A"
$ws.Range("D14").Value = "Jaylee, Cameron, Angelo, Ryan, Shane, Chandler"
$ws.Range("E14").Value = "Jaylee, Cameron, Ryan, Shane, Chandler, Angelo"
$ws.Range("B14:E14").WrapText = $true
$ws.Rows("14").RowHeight = 45

# --- Fix WEEK 12 / WEEK 13 rows (B11, B12, C12) ---
$ws.Range("B11").Value = "Shane, Jaylee, Cameron, Ryan, Angelo, Chandler"
$ws.Range("B12").Value = "Shane, Jaylee, Cameron, Ryan, Angelo, Chandler"
$ws.Range("C12").Value = "Jaylee, Cameron, Angelo, Ryan, Shane"

# --- Fix a couple of existing cells (E8, E10) that were mistakenly "N/A" ---
$ws.Range("E8").Value = "Shane, Jaylee, Cameron, Angelo, Ryan, Chandler"
$ws.Range("E10").Value = "Shane, Jaylee, Cameron, Angelo, Ryan, Chandler"

# --- Fill in WEEK 16 (row 15) ---
$ws.Range("A15").Value = "WEEK 16"
$ws.Range("B15").Value = "Shane, Jaylee, Cameron, Ryan, Angelo, Chandler"
$ws.Range("C15").Value = "N/A"
$ws.Range("D15").Value = "Jaylee, Cameron, Angelo, Ryan, Shane, Chandler"
$ws.Range("B15:E15").WrapText = $true
$ws.Rows("15").RowHeight = 45

# --- Update the view: scroll so row 3 is at the top, select I10 ---
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("I10").Select()
